# Apply updated crypto price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.183.54'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.213.87'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.17'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.77'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0897'
$ws.Range("E10").Value = '  +2.41%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.543.79'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.43'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.11'
$ws.Range("E14").Value = '  +2.74%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.232.67'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.089.15'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0930'
$ws.Range("E19").Value = '  +5.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  +2.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.93'
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '243.10'
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  +2.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +2.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.56'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.01'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.29'
$ws.Range("E29").Value = '  +2.03%  '
$ws.Range("E30").Value = '  +2.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0646'
$ws.Range("E35").Value = '  +4.08%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.55'
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.29'
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("E39").Value = '  +5.30%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("E41").Value = '  -3.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.53'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.20'
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0952'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.79'
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.34'
$ws.Range("E46").Value = '  -10.41%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.454.06'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.05'
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("E51").Value = '  +1.24%  '
